# Update the "Price" (column D) values on Sheet1 to reflect the latest
# symbol-list refresh, as produced by the scheduled GitHub Actions job.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of cell address -> new price text. Values are kept as text (not
# numbers) since the source cells are stored as inline strings.
$updates = [ordered]@{
    "D2"  = "268.53"
    "D3"  = "21.47"
    "D4"  = "6.192"
    "D5"  = "0.06160"
    "D6"  = "3.565"
    "D7"  = "6.539"
    "D8"  = "1.379"
    "D9"  = "0.8237"
    "D10" = "0.01351"
    "D12" = "0.08097"
    "D13" = "0.03354"
    "D14" = "0.03195"
    "D15" = "0.09220"
    "D16" = "3.739"
    "D17" = "0.001643"
    "D18" = "0.04670"
    "D19" = "0.006403"
    "D20" = "0.006192"
    "D21" = "0.001066"
    "D22" = "0.0001495"
    "D23" = "3.730"
    "D24" = "2.234"
    "D25" = "0.3305"
    "D26" = "0.1238"
    "D28" = "0.0002703"
    "D40" = "0.04649"
    "D41" = "0.006982"
    "D42" = "0.003987"
    "D43" = "0.1126"
    "D44" = "0.01177"
    "D45" = "0.00006041"
    "D46" = "0.0009865"
    "D48" = "0.7793"
    "D49" = "0.002175"
    "D50" = "0.00001894"
    "D51" = "0.01236"
}

foreach ($addr in $updates.Keys) {
    # Prefix with an apostrophe so Excel stores the value as literal text
    # (matching the source file's inline-string cells) instead of
    # re-interpreting the numeric-looking text as a number.
    $ws.Range($addr).Value = "'" + $updates[$addr]
}
